$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 464-476 ---
$pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-01-30", "18:12:18", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:20", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:23", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:28", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:33", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:38", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:43", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:48", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:53", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:12:58", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:13:04", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:13:09", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:13:14", "18:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 464
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]
    # Column A holds a date-shaped string ("YYYY-MM-DD"); force text so it
    # is stored as the literal string rather than an auto-converted date
    # serial. The other columns (time-of-day text, plain words) already
    # round-trip as literal text without help.
    $pir.Cells.Item($r, 1).NumberFormat = "@"
    $pir.Cells.Item($r, 1).Value = $row[0]
    $pir.Cells.Item($r, 2).Value = $row[1]
    $pir.Cells.Item($r, 3).Value = $row[2]
    $pir.Cells.Item($r, 4).Value = $row[3]
    $pir.Cells.Item($r, 5).Value = $row[4]
    $pir.Cells.Item($r, 6).Value = $row[5]
}

# --- Humidity sheet: append rows 305-313 ---
$hum = $wb.Worksheets.Item("Humidity")
$humRows = @(
    @("2026-01-30", "18:12:19", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:12:21", "18:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-30", "18:12:29", "18:00", "Bathroom", "86.8%", "Active"),
    @("2026-01-30", "18:12:39", "18:00", "Bathroom", "86.8%", "Active"),
    @("2026-01-30", "18:12:49", "18:00", "Bathroom", "86.5%", "Active"),
    @("2026-01-30", "18:12:59", "18:00", "Bathroom", "86.3%", "Active"),
    @("2026-01-30", "18:13:04", "18:00", "Bathroom", "86.2%", "Active"),
    @("2026-01-30", "18:13:10", "18:00", "Bathroom", "86.2%", "Active"),
    @("2026-01-30", "18:13:15", "18:00", "Bathroom", "86.2%", "Active")
)

$startRow = 305
for ($i = 0; $i -lt $humRows.Count; $i++) {
    $r = $startRow + $i
    $row = $humRows[$i]
    # Columns A ("YYYY-MM-DD") and E ("NN.N%") look numeric/date-shaped to
    # Excel's input parser, so force text on just those two cells before
    # assigning; B/C/D/F already round-trip as literal text untouched.
    $hum.Cells.Item($r, 1).NumberFormat = "@"
    $hum.Cells.Item($r, 5).NumberFormat = "@"
    $hum.Cells.Item($r, 1).Value = $row[0]
    $hum.Cells.Item($r, 2).Value = $row[1]
    $hum.Cells.Item($r, 3).Value = $row[2]
    $hum.Cells.Item($r, 4).Value = $row[3]
    $hum.Cells.Item($r, 5).Value = $row[4]
    $hum.Cells.Item($r, 6).Value = $row[5]
}
